$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same table of data and both
# need the "想去人数" (number of people interested) column (F) updated for
# rows 2, 3 and 5.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1763
    $ws.Range("F3").Value = 8057
    $ws.Range("F5").Value = 286
}
